# The worksheet is protected; unprotect it for the duration of the edits
# and re-protect it afterwards so the sheet behaves the same as before.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (A16).
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-09 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) columns for rows 2-13.
$ws.Range("D2").Value = 0.02722061445306007
$ws.Range("E2").Value = 0.01717433570965676

$ws.Range("D3").Value = 0.02154250301306366
$ws.Range("E3").Value = 0.002153625269203152

$ws.Range("D4").Value = 0.05688169304794199
$ws.Range("E4").Value = 0.005655042412818334

$ws.Range("D5").Value = 0.1403862563745155
$ws.Range("E5").Value = 0.004186046511627906

$ws.Range("D6").Value = 0.02014230618591906
$ws.Range("E6").Value = 0.02115655853314524

$ws.Range("D7").Value = 0.1292502067990937
$ws.Range("E7").Value = 0.01186399217221124

$ws.Range("D8").Value = 0.08868149984071229
$ws.Range("E8").Value = 0.01599702380952372

$ws.Range("D9").Value = 0.02921408783815131
$ws.Range("E9").Value = 0.02195871761089152

$ws.Range("D10").Value = 0.1013764039055216
$ws.Range("E10").Value = 0.02811639795517107

$ws.Range("D11").Value = 0.2960620710014231
$ws.Range("E11").Value = 0.009571788413098359

$ws.Range("D12").Value = 0.08924235754059776
$ws.Range("E12").Value = 0.01182994454713482

$ws.Range("D13").Value = 0.9999999999999999
$ws.Range("E13").Value = 0.01218284648341883

$ws.Protect()
